$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at N, shifting existing N:V columns to O:W
$ws.Columns("N:N").Insert()

# Header for the newly inserted column
$ws.Range("N1").Value2 = "Población activa, total"

# Populate the new column values (row -> value)
$ws.Range("N2").Value2 = 146729576
$ws.Range("N3").Value2 = 147698388
$ws.Range("N4").Value2 = 148526037
$ws.Range("N5").Value2 = 149140813
$ws.Range("N6").Value2 = 150211571
$ws.Range("N7").Value2 = 152071960
$ws.Range("N8").Value2 = 153931108
$ws.Range("N9").Value2 = 155224880
$ws.Range("N10").Value2 = 157006864
$ws.Range("N11").Value2 = 157109264
$ws.Range("N12").Value2 = 156903447
$ws.Range("N13").Value2 = 156975333
$ws.Range("N14").Value2 = 158636184
$ws.Range("N15").Value2 = 158755710
$ws.Range("N16").Value2 = 159532062
$ws.Range("N17").Value2 = 160596354
$ws.Range("N18").Value2 = 162547573
$ws.Range("N19").Value2 = 164268059
$ws.Range("N20").Value2 = 165483017
$ws.Range("N21").Value2 = 67652881
$ws.Range("N22").Value2 = 67601019
$ws.Range("N23").Value2 = 67173003
$ws.Range("N24").Value2 = 67065336
$ws.Range("N25").Value2 = 66725257
$ws.Range("N26").Value2 = 66734012
$ws.Range("N27").Value2 = 66802058
$ws.Range("N28").Value2 = 67143470
$ws.Range("N29").Value2 = 67100796
$ws.Range("N30").Value2 = 66914282
$ws.Range("N31").Value2 = 66663144
$ws.Range("N32").Value2 = 66043052
$ws.Range("N33").Value2 = 65639408
$ws.Range("N34").Value2 = 65970039
$ws.Range("N35").Value2 = 66175617
$ws.Range("N36").Value2 = 66363530
$ws.Range("N37").Value2 = 66850743
$ws.Range("N38").Value2 = 67288388
$ws.Range("N39").Value2 = 68358370
$ws.Range("N40").Value2 = 2760410612
$ws.Range("N41").Value2 = 2797040580
$ws.Range("N42").Value2 = 2838912969
$ws.Range("N43").Value2 = 2883821977
$ws.Range("N44").Value2 = 2932247001
$ws.Range("N45").Value2 = 2978118240
$ws.Range("N46").Value2 = 3014743223
$ws.Range("N47").Value2 = 3054910370
$ws.Range("N48").Value2 = 3090361014
$ws.Range("N49").Value2 = 3123387038
$ws.Range("N50").Value2 = 3152655443
$ws.Range("N51").Value2 = 3183946265
$ws.Range("N52").Value2 = 3218194184
$ws.Range("N53").Value2 = 3250599724
$ws.Range("N54").Value2 = 3283330912
$ws.Range("N55").Value2 = 3319710285
$ws.Range("N56").Value2 = 3354403108
$ws.Range("N57").Value2 = 3391647862
$ws.Range("N58").Value2 = 3427481440
$ws.Range("N59").Value2 = 199217046
$ws.Range("N60").Value2 = 199487250
$ws.Range("N61").Value2 = 200144636
$ws.Range("N62").Value2 = 201933394
$ws.Range("N63").Value2 = 203395227
$ws.Range("N64").Value2 = 205365646
$ws.Range("N65").Value2 = 206996811
$ws.Range("N66").Value2 = 208540629
$ws.Range("N67").Value2 = 210215521
$ws.Range("N68").Value2 = 210743587
$ws.Range("N69").Value2 = 211354253
$ws.Range("N70").Value2 = 210739030
$ws.Range("N71").Value2 = 212039861
$ws.Range("N72").Value2 = 212597037
$ws.Range("N73").Value2 = 213265285
$ws.Range("N74").Value2 = 213640167
$ws.Range("N75").Value2 = 214618312
$ws.Range("N76").Value2 = 215671511
$ws.Range("N77").Value2 = 216488451

Write-Output "done"
